$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "40-1="
$t.Cell(1,2).Range.Text = "77-29="
$t.Cell(1,3).Range.Text = "49+18="
$t.Cell(1,4).Range.Text = "45-6="
$t.Cell(1,5).Range.Text = "35+36="
$t.Cell(2,1).Range.Text = "55-7="
$t.Cell(2,2).Range.Text = "92-64="
$t.Cell(2,3).Range.Text = "94-8="
$t.Cell(2,4).Range.Text = "18+5="
$t.Cell(2,5).Range.Text = "40-17="
$t.Cell(3,1).Range.Text = "51-25="
$t.Cell(3,2).Range.Text = "92-18="
$t.Cell(3,3).Range.Text = "80-24="
$t.Cell(3,4).Range.Text = "33-4="
$t.Cell(3,5).Range.Text = "25+48="
$t.Cell(4,1).Range.Text = "30-2="
$t.Cell(4,2).Range.Text = "26+18="
$t.Cell(4,3).Range.Text = "60-41="
$t.Cell(4,4).Range.Text = "24+47="
$t.Cell(4,5).Range.Text = "73-19="
$t.Cell(5,1).Range.Text = "83+8="
$t.Cell(5,2).Range.Text = "64-57="
$t.Cell(5,3).Range.Text = "47+4="
$t.Cell(5,4).Range.Text = "53-45="
$t.Cell(5,5).Range.Text = "9+59="
$t.Cell(6,1).Range.Text = "57-18="
$t.Cell(6,2).Range.Text = "63-17="
$t.Cell(6,3).Range.Text = "90-35="
$t.Cell(6,4).Range.Text = "4+18="
$t.Cell(6,5).Range.Text = "33-24="
$t.Cell(7,1).Range.Text = "6+46="
$t.Cell(7,2).Range.Text = "51-8="
$t.Cell(7,3).Range.Text = "9+57="
$t.Cell(7,4).Range.Text = "91-75="
$t.Cell(7,5).Range.Text = "19+57="
$t.Cell(8,1).Range.Text = "65-8="
$t.Cell(8,2).Range.Text = "15+28="
$t.Cell(8,3).Range.Text = "53-27="
$t.Cell(8,4).Range.Text = "49+42="
$t.Cell(8,5).Range.Text = "75-67="
$t.Cell(9,1).Range.Text = "57+5="
$t.Cell(9,2).Range.Text = "94-68="
$t.Cell(9,3).Range.Text = "39+9="
$t.Cell(9,4).Range.Text = "23-9="
$t.Cell(9,5).Range.Text = "38+49="
$t.Cell(10,1).Range.Text = "4+89="
$t.Cell(10,2).Range.Text = "21-5="
$t.Cell(10,3).Range.Text = "82-38="
$t.Cell(10,4).Range.Text = "18+35="
$t.Cell(10,5).Range.Text = "37+48="
$t.Cell(11,1).Range.Text = "9+76="
$t.Cell(11,2).Range.Text = "51-2="
$t.Cell(11,3).Range.Text = "19+33="
$t.Cell(11,4).Range.Text = "96-7="
$t.Cell(11,5).Range.Text = "5+88="
$t.Cell(12,1).Range.Text = "51-43="
$t.Cell(12,2).Range.Text = "69+18="
$t.Cell(12,3).Range.Text = "39+23="
$t.Cell(12,4).Range.Text = "17+29="
$t.Cell(12,5).Range.Text = "55+9="
$t.Cell(13,1).Range.Text = "33+9="
$t.Cell(13,2).Range.Text = "80-57="
$t.Cell(13,3).Range.Text = "28+33="
$t.Cell(13,4).Range.Text = "18+45="
$t.Cell(13,5).Range.Text = "7+16="
$t.Cell(14,1).Range.Text = "51-43="
$t.Cell(14,2).Range.Text = "34-9="
$t.Cell(14,3).Range.Text = "65+16="
$t.Cell(14,4).Range.Text = "13-5="
$t.Cell(14,5).Range.Text = "12-5="
$t.Cell(15,1).Range.Text = "69+12="
$t.Cell(15,2).Range.Text = "73+18="
$t.Cell(15,3).Range.Text = "20-5="
$t.Cell(15,4).Range.Text = "5+57="
$t.Cell(15,5).Range.Text = "29+39="
$t.Cell(16,1).Range.Text = "73-25="
$t.Cell(16,2).Range.Text = "42-28="
$t.Cell(16,3).Range.Text = "33+39="
$t.Cell(16,4).Range.Text = "94-57="
$t.Cell(16,5).Range.Text = "58+27="
$t.Cell(17,1).Range.Text = "80-41="
$t.Cell(17,2).Range.Text = "57+24="
$t.Cell(17,3).Range.Text = "6+7="
$t.Cell(17,4).Range.Text = "52-39="
$t.Cell(17,5).Range.Text = "42-23="
$t.Cell(18,1).Range.Text = "87-18="
$t.Cell(18,2).Range.Text = "94-79="
$t.Cell(18,3).Range.Text = "53-45="
$t.Cell(18,4).Range.Text = "28+14="
$t.Cell(18,5).Range.Text = "81-12="
$t.Cell(19,1).Range.Text = "66-37="
$t.Cell(19,2).Range.Text = "25+16="
$t.Cell(19,3).Range.Text = "70-65="
$t.Cell(19,4).Range.Text = "18+78="
$t.Cell(19,5).Range.Text = "24+8="
$t.Cell(20,1).Range.Text = "47-8="
$t.Cell(20,2).Range.Text = "75-38="
$t.Cell(20,3).Range.Text = "61-37="
$t.Cell(20,4).Range.Text = "69+14="
$t.Cell(20,5).Range.Text = "44-36="
